# Update logInData.xlsx: the "email" column on Arkusz1 (A2:A4) now holds
# registration-form addresses (tst11/22/33 instead of tst111/222/333) and
# each one is turned into a clickable mailto: hyperlink, matching Excel's
# default "Hyperlink" cell style (underline + theme color).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New email addresses used for the "register" test-data rows.
$ws.Range("A2").Value = "tst11@test.com"
$ws.Range("A3").Value = "tst22@test.com"
$ws.Range("A4").Value = "tst33@test.com"

# Turn each address into a real mailto hyperlink (this also applies the
# built-in Hyperlink style to the cell).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:tst11@test.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:tst22@test.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:tst33@test.com")

# Leave the selection where the author ended up after the edit.
[void]$ws.Range("D5").Select()
